$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; existing rows 4-12 shift down to 5-13.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly price record.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44901
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103003
$ws.Range("J4").Value = "Damasco"
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = "$/bandeja 18 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 972
$ws.Range("T4").Value = 18
